$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 481, shifting existing rows 481:593 down to 482:594.
$ws.Rows.Item(481).Insert()

# Populate the newly inserted row 481 with the weekly price record
# (a copy of the original row 481 entry, updated with a new sample date,
# quality "Primera" and volume 500).
$ws.Range("A481").Value2 = 4
$ws.Range("B481").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C481").Value2 = "Los Lagos"
$ws.Range("D481").Value2 = 45173
$ws.Range("E481").Value2 = 10
$ws.Range("F481").Value2 = 100112008
$ws.Range("G481").Value2 = "Coliflor"
$ws.Range("H481").Value2 = "Sin especificar"
$ws.Range("I481").Value2 = "Primera"
$ws.Range("J481").Value2 = 500
$ws.Range("K481").Value2 = 1500
$ws.Range("L481").Value2 = 1500
$ws.Range("M481").Value2 = 1500
$ws.Range("N481").Value2 = "`$/unidad"
$ws.Range("O481").Value2 = "Región Metropolitana"
$ws.Range("P481").Value2 = 1500
$ws.Range("Q481").Value2 = 1
$ws.Range("R481").Value2 = "Hortaliza"
